# Apply attendance-count updates to Sheet1.
# For each listed row, set the indicated cell(s) from 0 to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    3  = @("G", "H")
    4  = @("D", "E")
    5  = @("H")
    6  = @("H")
    7  = @("H")
    8  = @("H")
    9  = @("D", "E")
    10 = @("D", "E")
    11 = @("D", "E")
    12 = @("H")
    13 = @("H")
    14 = @("H")
    15 = @("H")
    16 = @("H")
    17 = @("D", "E")
    18 = @("H")
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row]) {
        $ws.Range("$col$row").Value = 1
    }
}
